$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last-point fetch window moved: the oldest quarter (20180116) drops off
# the front, and the most recent several points (from 20191023 through 20200624)
# are trimmed off the tail, shifting every remaining row up by one.

# Rewrite A7:B37 with the shifted date/value pairs (rows 1-6 are unchanged)
$ws.Cells.Item(7, 1).Value = "'20180207"
$ws.Cells.Item(7, 2).Value = 1376.0
$ws.Cells.Item(8, 1).Value = "'20180312"
$ws.Cells.Item(8, 2).Value = 1718.0
$ws.Cells.Item(9, 1).Value = "'20180326"
$ws.Cells.Item(9, 2).Value = 1431.0
$ws.Cells.Item(10, 1).Value = "'20180413"
$ws.Cells.Item(10, 2).Value = 2022.0
$ws.Cells.Item(11, 1).Value = "'20180423"
$ws.Cells.Item(11, 2).Value = 1671.0
$ws.Cells.Item(12, 1).Value = "'20180528"
$ws.Cells.Item(12, 2).Value = 2289.0
$ws.Cells.Item(13, 1).Value = "'20180706"
$ws.Cells.Item(13, 2).Value = 1482.0
$ws.Cells.Item(14, 1).Value = "'20180713"
$ws.Cells.Item(14, 2).Value = 1737.0
$ws.Cells.Item(15, 1).Value = "'20180824"
$ws.Cells.Item(15, 2).Value = 1383.0
$ws.Cells.Item(16, 1).Value = "'20180912"
$ws.Cells.Item(16, 2).Value = 1520.0
$ws.Cells.Item(17, 1).Value = "'20181012"
$ws.Cells.Item(17, 2).Value = 1124.0
$ws.Cells.Item(18, 1).Value = "'20181113"
$ws.Cells.Item(18, 2).Value = 1485.0
$ws.Cells.Item(19, 1).Value = "'20181126"
$ws.Cells.Item(19, 2).Value = 1323.0
$ws.Cells.Item(20, 1).Value = "'20181211"
$ws.Cells.Item(20, 2).Value = 1941.0
$ws.Cells.Item(21, 1).Value = "'20181228"
$ws.Cells.Item(21, 2).Value = 1622.0
$ws.Cells.Item(22, 1).Value = "'20190109"
$ws.Cells.Item(22, 2).Value = 2041.0
$ws.Cells.Item(23, 1).Value = "'20190130"
$ws.Cells.Item(23, 2).Value = 1642.0
$ws.Cells.Item(24, 1).Value = "'20190225"
$ws.Cells.Item(24, 2).Value = 1968.0
$ws.Cells.Item(25, 1).Value = "'20190311"
$ws.Cells.Item(25, 2).Value = 1692.0
$ws.Cells.Item(26, 1).Value = "'20190319"
$ws.Cells.Item(26, 2).Value = 2216.0
$ws.Cells.Item(27, 1).Value = "'20190606"
$ws.Cells.Item(27, 2).Value = 1430.0
$ws.Cells.Item(28, 1).Value = "'20190724"
$ws.Cells.Item(28, 2).Value = 1808.0
$ws.Cells.Item(29, 1).Value = "'20190812"
$ws.Cells.Item(29, 2).Value = 1383.0
$ws.Cells.Item(30, 1).Value = "'20191121"
$ws.Cells.Item(30, 2).Value = 3392.0
$ws.Cells.Item(31, 1).Value = "'20191230"
$ws.Cells.Item(31, 2).Value = 1680.0
$ws.Cells.Item(32, 1).Value = "'20200121"
$ws.Cells.Item(32, 2).Value = 1865.0
$ws.Cells.Item(33, 1).Value = "'20200204"
$ws.Cells.Item(33, 2).Value = 1355.0
$ws.Cells.Item(34, 1).Value = "'20200225"
$ws.Cells.Item(34, 2).Value = 1978.0
$ws.Cells.Item(35, 1).Value = "'20200319"
$ws.Cells.Item(35, 2).Value = 1451.0
$ws.Cells.Item(36, 1).Value = "'20200521"
$ws.Cells.Item(36, 2).Value = 2700.0
$ws.Cells.Item(37, 1).Value = "'20200529"
$ws.Cells.Item(37, 2).Value = 2015.0

# Clear the now-unused trailing rows (previously rows 38-44)
$ws.Range("A38:B44").ClearContents()
